$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change -------------------------------------------------
# The sheet gains one row: a new data-only row is inserted at row 13
# holding the "Docentes responsaveis:" value (previously mis-aligned one
# row too high, duplicating the "Objetivos:" value). Inserting a row here
# shifts every row from the old 13 down to 14, which automatically carries
# the previously-correct row heights along for rows 14 and below.
$ws.Rows.Item(13).Insert()

# The insert copies formatting down from row 12 (label-only row), which
# leaves a phantom styled-but-empty A13 cell; the target layout has no
# cell in column A on this row at all.
$ws.Range("A13").Clear()

# Populate the newly inserted row 13 with the "Docentes responsaveis:"
# value (Portuguese/English columns share the same text).
$ws.Range("B13").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C13").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("B13").WrapText = $true
$ws.Range("B13").Font.Bold = $false

# --- Content fixes ------------------------------------------------------

# Objetivos: value was wrongly the professor's name; replace with the
# actual objectives text.
$objetivos = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes relacionados as fronteira da engenharia química e a suas interfaces com outras ciências."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Programa resumido: / Short syllabus value was wrongly "Semestral";
# Programa: / Syllabus value was wrongly the activation date. Both should
# hold the short topic description.
$topicos = "Tópicos atuais e relevantes relacionados as fronteira da engenharia química e a suas interfaces com outras ciências."
$ws.Range("B14").Value = $topicos
$ws.Range("C14").Value = $topicos
$ws.Range("B16").Value = $topicos
$ws.Range("C16").Value = $topicos

# Metodo: value was wrongly the professor's name; it should hold the
# teaching-method description (which used to sit one row down, on
# "Criterio:").
$metodo = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Criterio: value should be "Provas e trabalhos." (previously one row off,
# sitting on "Norma de recuperacao:").
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# Norma de recuperacao: value should be the single-exam passing rule
# (previously one row off, sitting on "Bibliografia:").
$recuperacao = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao

# Bibliografia: gains real bibliography text (previously empty/misaligned).
$biblio = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas de Engenharia Química."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
